# [FEATURE] add fiori reuse component for UploadCollection
#
# This script edits the TestDataSet / TechSet workbook to:
#  - rename the "Guid" column to "RowId" (TestDataSet) and replace the GUID
#    values with simple zero-padded row numbers ("01".."10")
#  - replace the "Status" column's GUID-like values with simple zero-padded
#    codes ("00","01","02") and drop the old text-format styling on that
#    column (TestDataSet + TechSet)
#  - turn TechSet's HasApprove/HasChange/HasReplace columns into real
#    booleans
#  - add a brand-new "FileSet" worksheet + table with upload-collection
#    sample data (RowId, FileId, FileName, UserName, FileSize)
#  - tidy up the sheet selections/active tab so FileSet ends up selected

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) TestDataSet (sheet 1) -------------------------------------------
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TestDataSet")

$ws1.Range("A1").Value = "RowId"

$rowIds = @("01","02","03","04","05","06","07","08","09","10")
$strings = @(
    "Reprehenderit sint velit. Ipsum laboris.",
    "Lorem aliquip id pariatur.",
    "Mollit excepteur.",
    "Laboris nulla sunt duis. Aliqua anim.",
    "Labore qui id. Nostrud laborum officia.",
    "Ullamco est. Velit ullamco eiusmod ad.",
    "Labore aliqua.",
    "Ex velit deserunt. Minim nulla sit et.",
    "Eu dolor. Dolor laboris.",
    "Amet velit deserunt."
)
$statuses = @("00","01","02","00","01","02","00","01","02","02")

# force text storage (values are zero-padded codes, not numbers) while
# writing, then strip the formatting back off so the cells end up with
# the workbook's default (unstyled) look, matching the target cells
$ws1.Range("A2:A11").NumberFormat = "@"
$ws1.Range("C2:C11").NumberFormat = "@"

for ($i = 0; $i -lt 10; $i++) {
    $r = $i + 2
    $ws1.Range("A$r").Value = $rowIds[$i]
    $ws1.Range("B$r").Value = $strings[$i]
    $ws1.Range("C$r").Value = $statuses[$i]
}

# Status (and RowId) columns no longer carry the explicit "@" text
# numberformat style that used to be applied via conditional formatting
$ws1.Range("A2:A11").ClearFormats()
$ws1.Range("C2:C11").ClearFormats()

# Drop the stale selection left over from editing
$ws1.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) TechSet (sheet 2) -------------------------------------------------
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("TechSet")

$ws2.Range("A2:A4").NumberFormat = "@"
$ws2.Range("A2").Value = "00"
$ws2.Range("A3").Value = "01"
$ws2.Range("A4").Value = "02"
$ws2.Range("A2:A4").ClearFormats()

$ws2.Range("B2").Value = $true
$ws2.Range("C2").Value = $false
$ws2.Range("D2").Value = $false

$ws2.Range("B3").Value = $false
$ws2.Range("C3").Value = $true
$ws2.Range("D3").Value = $false

$ws2.Range("B4").Value = $false
$ws2.Range("C4").Value = $false
$ws2.Range("D4").Value = $true

$ws2.Range("A1").Select()

# ---------------------------------------------------------------------
# 3) New FileSet worksheet + table ------------------------------------
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws3.Name = "FileSet"

$ws3.Range("A1").Value = "RowId"
$ws3.Range("B1").Value = "FileId"
$ws3.Range("C1").Value = "FileName"
$ws3.Range("D1").Value = "UserName"
$ws3.Range("E1").Value = "FileSize"

$fileRows = @(
    @("01","001","File 1","Admin","100Kb"),
    @("01","002","File 2","User","1M"),
    @("01","003","File 3","Vendor","100Kb"),
    @("01","004","File 4","Guest","100Kb")
)

$ws3.Range("A2:B5").NumberFormat = "@"

for ($i = 0; $i -lt $fileRows.Count; $i++) {
    $r = $i + 2
    $row = $fileRows[$i]
    $ws3.Range("A$r").Value = $row[0]
    $ws3.Range("B$r").Value = $row[1]
    $ws3.Range("C$r").Value = $row[2]
    $ws3.Range("D$r").Value = $row[3]
    $ws3.Range("E$r").Value = $row[4]
}

$ws3.Range("A2:B5").ClearFormats()

$fileTbl = $ws3.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws3.Range("A1:E5"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$fileTbl.Name = "FileSet"

$ws3.Range("E8").Select()
$ws3.Activate()
